$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text/"@" number format on price cells whose new values would
# otherwise be auto-parsed as numbers by Excel, losing exact text formatting
# (e.g. trailing zeros, trailing-dot grouped numbers).
$textCells = @("D4","D5","D7","D8","D9","D10","D11","D13","D14","D15","D16","D17","D18","D19","D23","D24","D25","D27","D28","D29","D30","D31","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell, in sheet (row, column) order.
$ws.Range("D2").Value = '26.949.70'
$ws.Range("E2").Value = '  +0.55%  '
$ws.Range("D3").Value = '1.876.14'
$ws.Range("E3").Value = '  +1.04%  '
$ws.Range("D4").Value = '0.9990'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '306.38'
$ws.Range("E5").Value = '  +0.44%  '
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("D7").Value = '0.5155'
$ws.Range("E7").Value = '  +1.28%  '
$ws.Range("D8").Value = '0.3715'
$ws.Range("E8").Value = '  +1.66%  '
$ws.Range("D9").Value = '0.07189'
$ws.Range("E9").Value = '  +0.77%  '
$ws.Range("D10").Value = '0.8985'
$ws.Range("E10").Value = '  +1.17%  '
$ws.Range("D11").Value = '20.69'
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.897.43'
$ws.Range("E12").Value = '  +2.19%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.07561'
$ws.Range("E13").Value = '  +0.78%  '
$ws.Range("D14").Value = '94.91'
$ws.Range("E14").Value = '  +4.52%  '
$ws.Range("D15").Value = '5.249'
$ws.Range("E15").Value = '  +0.13%  '
$ws.Range("D16").Value = '0.9996'
$ws.Range("D17").Value = '0.000008481'
$ws.Range("E17").Value = '  -0.66%  '
$ws.Range("D18").Value = '14.24'
$ws.Range("E18").Value = '  +1.47%  '
$ws.Range("D19").Value = '0.9993'
$ws.Range("E19").Value = '  -0.14%  '
$ws.Range("D20").Value = '26.969.55'
$ws.Range("E20").Value = '  +0.46%  '
$ws.Range("E21").Value = '  +0.49%  '
$ws.Range("D22").Value = '2.122.51'
$ws.Range("E22").Value = '  +1.29%  '
$ws.Range("D23").Value = '10.40'
$ws.Range("E23").Value = '  +1.44%  '
$ws.Range("D24").Value = '6.439'
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").Value = '146.05'
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("D27").Value = '18.03'
$ws.Range("E27").Value = '  +1.11%  '
$ws.Range("D28").Value = '2.107'
$ws.Range("E28").Value = '  +3.29%  '
$ws.Range("D29").Value = '114.51'
$ws.Range("E29").Value = '  +1.47%  '
$ws.Range("D30").Value = '4.897'
$ws.Range("E30").Value = '  +4.92%  '
$ws.Range("D31").Value = '4.753'
$ws.Range("E31").Value = '  +2.89%  '
$ws.Range("E32").Value = '  -0.73%  '
$ws.Range("D33").Value = '0.05034'
$ws.Range("E33").Value = '  -1.41%  '
$ws.Range("D34").Value = '0.7540'
$ws.Range("E34").Value = '  +3.06%  '
$ws.Range("D35").Value = '3.000'
$ws.Range("E35").Value = '  -2.26%  '
$ws.Range("D36").Value = '1.172'
$ws.Range("E36").Value = '  +1.92%  '
$ws.Range("D37").Value = '3.275'
$ws.Range("E37").Value = '  +2.92%  '
$ws.Range("D38").Value = '0.01995'
$ws.Range("E38").Value = '  -0.96%  '
$ws.Range("D39").Value = '0.5580'
$ws.Range("E39").Value = '  +5.54%  '
$ws.Range("D40").Value = '2.487'
$ws.Range("E40").Value = '  +1.27%  '
$ws.Range("D41").Value = '1.074'
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").Value = '6.571'
$ws.Range("E42").Value = '  +1.80%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '116.39'
$ws.Range("E43").Value = '  -1.22%  '
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").Value = '8.733'
$ws.Range("E44").Value = '  +4.55%  '
$ws.Range("E45").Value = '  +1.96%  '
$ws.Range("D46").Value = '0.4771'
$ws.Range("E46").Value = '  +2.68%  '
$ws.Range("D47").Value = '0.9992'
$ws.Range("E47").Value = '  -0.11%  '
$ws.Range("D48").Value = '10.14'
$ws.Range("E48").Value = '  +2.21%  '
$ws.Range("D49").Value = '1.563'
$ws.Range("D50").Value = '37.13'
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("D51").Value = '63.38'
$ws.Range("E51").Value = '  +0.46%  '
